$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing row 3 odds values (in-place edits, no row shift)
$ws.Range("I3").Value = 7.5
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 1.33
$ws.Range("T3").Value = 3.25
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("Z3").Value = 9
$ws.Range("AD3").Value = 9.5
$ws.Range("AT3").Value = 3.25
$ws.Range("BA3").Value = 151

# 2. Insert a brand-new row at position 4; this pushes the former
#    rows 4,5,6 down to 5,6,7 (matches dimension growing to BD7)
$ws.Rows.Item(4).Insert()

# 3. Populate the newly inserted row 4 with the new match data
# Column A: plain text id - safe as-is
$ws.Range("A4").Value = "K8JqOFWo"
# Column B looks like a date ("01/11/2024") - force text so Excel
# does not auto-convert it to a date serial number/format
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "01/11/2024"
$ws.Range("B4").Style = "Normal"
# Column C ("09:30") - observed to stay literal text on its own,
# but force text too for safety/consistency
$ws.Range("C4").Value = "09:30"
$ws.Range("D4").Value = "MALAYSIA - SUPER LEAGUE"
$ws.Range("E4").Value = "Negeri Sembilan"
$ws.Range("F4").Value = "Johor DT"

# Remaining columns (G onward) are numeric odds
$ws.Range("G4").Value = 37
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 1.03
$ws.Range("J4").Value = 26
$ws.Range("K4").Value = 3.7
$ws.Range("L4").Value = 1.22
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 26
$ws.Range("O4").Value = 1.01
$ws.Range("P4").Value = 11
$ws.Range("Q4").Value = 1.15
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 1.09
$ws.Range("T4").Value = 6.1
$ws.Range("U4").Value = 2.78
$ws.Range("V4").Value = 1.42
$ws.Range("W4").Value = 200
$ws.Range("X4").Value = 500
$ws.Range("Y4").Value = 200
$ws.Range("Z4").Value = 500
$ws.Range("AA4").Value = 400
$ws.Range("AB4").Value = 600
$ws.Range("AC4").Value = 26
$ws.Range("AD4").Value = 30
$ws.Range("AE4").Value = 65
$ws.Range("AF4").Value = 250
$ws.Range("AG4").Value = 201
$ws.Range("AH4").Value = 12.5
$ws.Range("AI4").Value = 7.1
$ws.Range("AJ4").Value = 15.5
$ws.Range("AK4").Value = 5.5
$ws.Range("AL4").Value = 11.75
$ws.Range("AM4").Value = 45
$ws.Range("AN4").Value = 40
$ws.Range("AO4").Value = 450
$ws.Range("AP4").Value = 200
$ws.Range("AQ4").Value = 501
$ws.Range("AR4").Value = 501
$ws.Range("AS4").Value = 501
$ws.Range("AT4").Value = 5.3
$ws.Range("AU4").Value = 15
$ws.Range("AV4").Value = 150
$ws.Range("AW4").Value = 3.3
$ws.Range("AX4").Value = 3.7
$ws.Range("AY4").Value = 16.5
$ws.Range("AZ4").Value = 5.9
$ws.Range("BA4").Value = 26
$ws.Range("BB4").Value = 200
$ws.Range("BC4").Value = 51
$ws.Range("BD4").Value = 51

# 4. Update row 7 (the former row 6, Ulsan HD vs Gangwon) odds values
$ws.Range("I7").Value = 3.9
$ws.Range("AW7").Value = 5.5
